# Update "想去人数" (interested count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1102
$ws1.Range("F4").Value = 1757
$ws1.Range("F5").Value = 782
$ws1.Range("F6").Value = 171
$ws1.Range("F7").Value = 203

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1102
$ws4.Range("F4").Value = 1757
$ws4.Range("F6").Value = 782
$ws4.Range("F7").Value = 171
$ws4.Range("F8").Value = 203

$wb.Save()
